# Update countries & provincias Spain
# Applies the COVID-19 stats refresh: updates the "last updated" timestamp,
# re-ranks a handful of countries whose case counts now place them in a
# different position in the table (Portugal/Etiopia, Birmania/Gabon/Haiti,
# Islas Malvinas/Montserrat), and refreshes the numeric columns for every
# row whose totals changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 16:45"

$ws.Range("B4").Value = 7191643
$ws.Range("C4").Value = 6172
$ws.Range("D4").Value = 4440410
$ws.Range("E4").Value = 2543612
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 207621

$ws.Range("B5").Value = 5843349
$ws.Range("C5").Value = 27246
$ws.Range("D5").Value = 4779658
$ws.Range("E5").Value = 971104
$ws.Range("G5").Value = 270
$ws.Range("H5").Value = 92587

$ws.Range("B15").Value = 453868
$ws.Range("C15").Value = 2222
$ws.Range("D15").Value = 428580
$ws.Range("E15").Value = 12761
$ws.Range("G15").Value = 58
$ws.Range("H15").Value = 12527

$ws.Range("B25").Value = 282163
$ws.Range("C25").Value = 818
$ws.Range("E25").Value = 23141

$ws.Range("B27").Value = 217899
$ws.Range("C27").Value = 5784
$ws.Range("D27").Value = 153574
$ws.Range("E27").Value = 62913
$ws.Range("G27").Value = 34
$ws.Range("H27").Value = 1412

$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 72055
$ws.Range("C51").Value = 899
$ws.Range("D51").Value = 47003
$ws.Range("E51").Value = 23116
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 1936

$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 71687
$ws.Range("D52").Value = 29461
$ws.Range("E52").Value = 41078
$ws.Range("H52").Value = 1148

$ws.Range("E61").Value = 7200
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 2064

$ws.Range("B74").Value = 33238
$ws.Range("C74").Value = 75
$ws.Range("E74").Value = 956
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 746

$ws.Range("B86").Value = 17343
$ws.Range("C86").Value = 143
$ws.Range("D86").Value = 14422
$ws.Range("E86").Value = 2200
$ws.Range("G86").Value = 9
$ws.Range("H86").Value = 721

$ws.Range("B91").Value = 14515
$ws.Range("C91").Value = 24
$ws.Range("E91").Value = 540

$ws.Range("B93").Value = 13478
$ws.Range("C93").Value = 72
$ws.Range("E93").Value = 2837

$ws.Range("A105").Value = "Birmania"
$ws.Range("B105").Value = 9112
$ws.Range("C105").Value = 768
$ws.Range("D105").Value = 2381
$ws.Range("E105").Value = 6557
$ws.Range("G105").Value = 24
$ws.Range("H105").Value = 174

$ws.Range("A106").Value = "Gabon"
$ws.Range("B106").Value = 8716
$ws.Range("D106").Value = 7906
$ws.Range("E106").Value = 756
$ws.Range("H106").Value = 54

$ws.Range("A107").Value = "Haiti"
$ws.Range("B107").Value = 8684
$ws.Range("C107").Value = 16
$ws.Range("D107").Value = 6551
$ws.Range("E107").Value = 1906
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 227

$ws.Range("B116").Value = 5723
$ws.Range("C116").Value = 135
$ws.Range("D116").Value = 1496
$ws.Range("E116").Value = 4147
$ws.Range("G116").Value = 3
$ws.Range("H116").Value = 80

$ws.Range("B120").Value = 5350
$ws.Range("C120").Value = 40
$ws.Range("D120").Value = 4658
$ws.Range("E120").Value = 574

$ws.Range("B132").Value = 4260
$ws.Range("C132").Value = 25
$ws.Range("D132").Value = 2122
$ws.Range("E132").Value = 2069
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 69

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

